$d = $word.ActiveDocument

$d.Content.Find.Execute("54×27=", $true, $false, $false, $false, $false, $true, 1, $false, "29×19=", 2) | Out-Null
$d.Content.Find.Execute("25×91=", $true, $false, $false, $false, $false, $true, 1, $false, "61×60=", 2) | Out-Null
$d.Content.Find.Execute("57×19=", $true, $false, $false, $false, $false, $true, 1, $false, "68×65=", 2) | Out-Null
$d.Content.Find.Execute("84×37=", $true, $false, $false, $false, $false, $true, 1, $false, "23×99=", 2) | Out-Null
$d.Content.Find.Execute("68×64=", $true, $false, $false, $false, $false, $true, 1, $false, "73×89=", 2) | Out-Null
$d.Content.Find.Execute("90×28=", $true, $false, $false, $false, $false, $true, 1, $false, "45×27=", 2) | Out-Null
$d.Content.Find.Execute("33×11=", $true, $false, $false, $false, $false, $true, 1, $false, "27×67=", 2) | Out-Null
$d.Content.Find.Execute("84×83=", $true, $false, $false, $false, $false, $true, 1, $false, "65×99=", 2) | Out-Null
$d.Content.Find.Execute("98×64=", $true, $false, $false, $false, $false, $true, 1, $false, "30×84=", 2) | Out-Null
$d.Content.Find.Execute("77×73=", $true, $false, $false, $false, $false, $true, 1, $false, "54×21=", 2) | Out-Null
$d.Content.Find.Execute("72×29=", $true, $false, $false, $false, $false, $true, 1, $false, "45×60=", 2) | Out-Null
$d.Content.Find.Execute("17×91=", $true, $false, $false, $false, $false, $true, 1, $false, "84×53=", 2) | Out-Null
$d.Content.Find.Execute("37×12=", $true, $false, $false, $false, $false, $true, 1, $false, "34×29=", 2) | Out-Null
$d.Content.Find.Execute("63×58=", $true, $false, $false, $false, $false, $true, 1, $false, "34×91=", 2) | Out-Null
$d.Content.Find.Execute("51×69=", $true, $false, $false, $false, $false, $true, 1, $false, "28×12=", 2) | Out-Null
$d.Content.Find.Execute("29×96=", $true, $false, $false, $false, $false, $true, 1, $false, "57×95=", 2) | Out-Null
$d.Content.Find.Execute("21×40=", $true, $false, $false, $false, $false, $true, 1, $false, "98×27=", 2) | Out-Null
$d.Content.Find.Execute("79×20=", $true, $false, $false, $false, $false, $true, 1, $false, "20×89=", 2) | Out-Null
$d.Content.Find.Execute("16×61=", $true, $false, $false, $false, $false, $true, 1, $false, "85×55=", 2) | Out-Null
$d.Content.Find.Execute("49×63=", $true, $false, $false, $false, $false, $true, 1, $false, "31×52=", 2) | Out-Null
$d.Content.Find.Execute("83×63=", $true, $false, $false, $false, $false, $true, 1, $false, "83×73=", 2) | Out-Null
$d.Content.Find.Execute("26×28=", $true, $false, $false, $false, $false, $true, 1, $false, "60×94=", 2) | Out-Null
$d.Content.Find.Execute("46×90=", $true, $false, $false, $false, $false, $true, 1, $false, "51×25=", 2) | Out-Null
$d.Content.Find.Execute("21×99=", $true, $false, $false, $false, $false, $true, 1, $false, "44×20=", 2) | Out-Null
$d.Content.Find.Execute("31×99=", $true, $false, $false, $false, $false, $true, 1, $false, "58×22=", 2) | Out-Null
